$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) AddOpportunity sheet: duplicate row 2 into row 3 (keeping every
#    format/style) and rename the Client/Subject columns to the new
#    "Triangle Capital Corporation" company, leaving JobType as
#    "Sellside".
# ---------------------------------------------------------------------
$wsOpp = $wb.Worksheets.Item("AddOpportunity")
$wsOpp.Range("A2:AD2").Copy($wsOpp.Range("A3:AD3"))
$wsOpp.Range("A3").Value = "Triangle Capital Corporation"
$wsOpp.Range("B3").Value = "Triangle Capital Corporation"
$wsOpp.Columns.Item(1).ColumnWidth = 24.7

# ---------------------------------------------------------------------
# 2) AddCounterparty sheet: add the matching counterparty row for the
#    new Triangle Capital Corporation company.
# ---------------------------------------------------------------------
$wsCp = $wb.Worksheets.Item("AddCounterparty")
$wsCp.Range("A3").Value = "Triangle Capital Corporation"
$wsCp.Range("B3").Value = "Financial"
$wsCp.Columns.Item(1).ColumnWidth = 24.7
$wsCp.Range("E4").Select()

# ---------------------------------------------------------------------
# 3) FlagReason sheet: clear the two "capital provider / operating
#    company ... Brian Miller" comment cells (C2, D2) so those two
#    shared strings become unreferenced and drop out of the workbook.
# ---------------------------------------------------------------------
$wsFlag = $wb.Worksheets.Item("FlagReason")
$wsFlag.Range("C2").ClearContents()
$wsFlag.Range("D2").ClearContents()
$wsFlag.Rows.Item(2).RowHeight = 28.8
$wsFlag.Range("C11").Select()

# ---------------------------------------------------------------------
# 4) Warning sheet: replace the second comment cell with the new
#    "Subject is typically considered a potential round trip" message.
# ---------------------------------------------------------------------
$wsWarn = $wb.Worksheets.Item("Warning")
$wsWarn.Range("B2").Value = "A Subject is typically considered a potential round trip if it is an operating company acquired either by a Private Equity firm or by a PE-owned operating company. The Subject is not listed as an Operating Company. If you still want to consider them a round trip candidate no change is needed; otherwise, please change the selection."
$wsWarn.Rows.Item(2).RowHeight = 129.6
$wsWarn.Range("B5").Select()

# ---------------------------------------------------------------------
# 5) Switch the active tab from Bid to FlagReason, matching the new
#    workbook view state.
# ---------------------------------------------------------------------
$wsFlag.Activate()
